# Insert a new snapshot column right before the "nom" column (EQ),
# shifting "nom" (old EQ) to ER and "url_produit" (old ER) to ES.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("EQ").Insert()

# New header timestamp for the freshly inserted column.
$ws.Range("EQ1").Value2 = "2026-02-03 16:39:03"

# The new snapshot column repeats the latest known price (previous
# snapshot column, now at EP) for every product row that already had a
# price recorded there; rows where EP is blank stay blank.
$ws.Range("EQ2:EQ80").Value2 = $ws.Range("EP2:EP80").Value2
